$d = $word.ActiveDocument

# --- 1) Text edits -------------------------------------------------------
# "...that as you can see later..."      -> "...that as you can be seen later..."
$d.Content.Find.Execute("can see later", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "can be seen later", 2) | Out-Null

# "...has given us the best results compared..." -> "...best results, compared..."
$d.Content.Find.Execute("the best results compared", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "the best results, compared", 2) | Out-Null

# --- 2) Re-create the run boundaries seen in the authored edit ----------
# Word breaks a paragraph's single run into several runs at every point the
# cursor stopped during the edit session. We reproduce those breakpoints by
# dropping a temporary bookmark at each boundary (which forces a run split)
# and then removing the bookmark again (the split itself persists).
function Split-At([string]$beforeText) {
    $r = $d.Content
    $r.Find.Execute($beforeText) | Out-Null
    $pos = $r.End
    $markRange = $d.Range($pos, $pos)
    $d.Bookmarks.Add("_tmpSplit", $markRange) | Out-Null
    $d.Bookmarks.Item("_tmpSplit").Delete()
}

Split-At "we encountered"
Split-At "this algorithm, that as you can "
Split-At "this algorithm, that as you can be seen"
Split-At "has given us the best results"
Split-At "has given us the best results,"

# --- 3) Relocate the "_GoBack" bookmark ----------------------------------
# Word stamps "_GoBack" at the site of the most recent edit; here that is
# right before "gives us the best compression...". Adding a bookmark with
# this reserved name also removes the pre-existing "_GoBack" elsewhere in
# the document (bookmark names are unique), matching the paragraph at the
# end of the document losing its bookmark.
$target = $d.Content
$target.Find.Execute("gives us the best compression") | Out-Null
$goBackRange = $d.Range($target.Start, $target.Start)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
